$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("text_coercion")

# Add a new row of data: a numeric "student number" value paired with its text label.
$ws.Range("A9").Value = 36436153
$ws.Range("B9").Value = "student number"

# Move/collapse the active selection to A9, mirroring the recorded cursor position.
$ws.Range("A9").Select()
